# Slide 1: "TextBox 2" (shape id 3) holds the attendance-password text box.
# Its second paragraph is currently a blank-line placeholder ("__________");
# replace it with the actual password text "debug" while leaving the run's
# existing formatting (size, highlight, etc.) untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$passwordPara = $tr.Paragraphs(2, 1)
$passwordPara.Text = "debug"
